$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels: player_id -> player_id_x, birth_year -> birth_year_x
$ws.Cells.Item(1, 3).Value2 = "player_id_x"
$ws.Cells.Item(1, 5).Value2 = "birth_year_x"

# Update player_id values in column C (rows 2-23) from 3463 to 3462
$lastRow = $ws.Cells.Item(1, 1).End(4).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 3463) {
        $cell.Value2 = 3462
    }
}
